$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Contest 19 RCB vs DC results + new Contest 28 "RCB vs KKR" row
# ---------------------------------------------------------------------------

# 1) Insert a new blank row at row 37 (shifts old rows 37-45 down to 38-46).
#    Row 37 was the last row inside the SUM(D10:D37) ranges (and similar),
#    so inserting here makes Excel auto-extend those sums to D10:D38 etc,
#    matching the target workbook.
$ws.Rows("37:37").Insert()

# 2) Copy formats (not content) from row 38 (the shifted-down old template
#    row, which already has the correct borders/fills) into the newly
#    blank row 37 so it matches the rest of the results table.
$fmtCols = @("A","B","C","D","E","G","H","J","K","M","N","P","Q","S","T")
foreach ($col in $fmtCols) {
    $ws.Range($col + "38").Copy()
    $ws.Range($col + "37").PasteSpecial(-4122)
}

# 3) Fill in match 28 details on row 37 ("RCB vs KKR")
$ws.Range("A37").Value = 28
$ws.Range("B37").Value = 1
$ws.Range("C37").Value = "RCB vs KKR"

$pairs = @{ "D" = "E"; "G" = "H"; "J" = "K"; "M" = "N"; "P" = "Q"; "S" = "T" }
foreach ($col in $pairs.Keys) {
    $src = $pairs[$col]
    $formula = '=IF(ISERROR(VLOOKUP(RANK(' + $src + '37, ($T37,$Q37,$N37,$K37,$H37,$E37), 0),  score, 2, FALSE)),"",VLOOKUP(RANK(' + $src + '37, ($T37,$Q37,$N37,$K37,$H37,$E37), 0),  score, 2, FALSE))'
    $ws.Range($col + "37").Formula = $formula
}

# 4) Fill in the results for match 19 (row 28, "RCB vs DC")
$ws.Range("E28").Value = 0
$ws.Range("H28").Value = 20
$ws.Range("K28").Value = 80
$ws.Range("N28").Value = 40
$ws.Range("Q28").Value = 100
$ws.Range("T28").Value = 60

# 5) Fix up the active-cell selection (was U41, table grew by one row so the
#    grand-total cell is now U42).
$ws.Range("U42").Select()

# 6) Re-anchor the conditional formatting rules that highlight the six
#    "Total" cells -- they still point at row 41, move them down to row 42.
$totalCols = @("E", "H", "K", "N", "Q", "T")
foreach ($col in $totalCols) {
    $oldRange = $ws.Range($col + "41")
    $newRange = $ws.Range($col + "42")
    $fcs = $oldRange.FormatConditions()
    $count = $fcs.Count()
    for ($i = 1; $i -le $count; $i++) {
        $fc = $fcs.Item($i)
        $fc.ModifyAppliesToRange($newRange)
    }
}
